$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(4).Insert()
$ws.Range("D5:D102").Value2 = 1234567
$ws.Columns("D:K").AutoFit()
Write-Host ("After autofit - ColumnWidth D: " + $ws.Columns.Item(4).ColumnWidth)
Write-Host ("After autofit - ColumnWidth E: " + $ws.Columns.Item(5).ColumnWidth)
Write-Host ("After autofit - ColumnWidth B: " + $ws.Columns.Item(2).ColumnWidth)
Write-Host ("After autofit - ColumnWidth C: " + $ws.Columns.Item(3).ColumnWidth)
